{"js": "// Replace each two-digit-multiplication answer in the table with its\n// updated counterpart. Every \"NN\u00d7NN=NNNN\" string in the document is\n// unique, so an exact, case-sensitive whole-match search safely\n// identifies the single cell to update for each pair.\nconst pairs = [[\"91\u00d726=2366\", \"90\u00d748=4320\"], [\"14\u00d756=784\", \"23\u00d764=1472\"], [\"80\u00d732=2560\", \"46\u00d791=4186\"], [\"93\u00d735=3255\", \"24\u00d730=720\"], [\"67\u00d784=5628\", \"95\u00d750=4750\"], [\"65\u00d795=6175\", \"59\u00d743=2537\"], [\"56\u00d727=1512\", \"70\u00d743=3010\"], [\"24\u00d762=1488\", \"26\u00d799=2574\"], [\"42\u00d723=966\", \"63\u00d729=1827\"], [\"65\u00d766=4290\", \"41\u00d791=3731\"], [\"89\u00d786=7654\", \"81\u00d752=4212\"], [\"75\u00d793=6975\", \"28\u00d775=2100\"], [\"18\u00d758=1044\", \"97\u00d772=6984\"], [\"54\u00d728=1512\", \"51\u00d781=4131\"], [\"95\u00d758=5510\", \"44\u00d793=4092\"], [\"15\u00d711=165\", \"92\u00d732=2944\"], [\"40\u00d712=480\", \"81\u00d751=4131\"], [\"41\u00d764=2624\", \"71\u00d783=5893\"], [\"57\u00d722=1254\", \"46\u00d790=4140\"], [\"81\u00d794=7614\", \"82\u00d771=5822\"], [\"47\u00d766=3102\", \"70\u00d768=4760\"], [\"89\u00d749=4361\", \"28\u00d725=700\"], [\"80\u00d738=3040\", \"68\u00d721=1428\"], [\"68\u00d760=4080\", \"39\u00d716=624\"], [\"94\u00d738=3572\", \"55\u00d793=5115\"]];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of pairs) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Could not find text to replace: ${oldText}`);\n  }\n\n  for (const r of results.items) {\n    r.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each two-digit-multiplication answer in the table with its\n# updated counterpart. Every \"NN\u00d7NN=NNNN\" string in the document is\n# unique, so an exact, case-sensitive Find/Replace safely targets the\n# single cell that needs to change for each pair.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"91\u00d726=2366\", \"90\u00d748=4320\"),\n    @(\"14\u00d756=784\", \"23\u00d764=1472\"),\n    @(\"80\u00d732=2560\", \"46\u00d791=4186\"),\n    @(\"93\u00d735=3255\", \"24\u00d730=720\"),\n    @(\"67\u00d784=5628\", \"95\u00d750=4750\"),\n    @(\"65\u00d795=6175\", \"59\u00d743=2537\"),\n    @(\"56\u00d727=1512\", \"70\u00d743=3010\"),\n    @(\"24\u00d762=1488\", \"26\u00d799=2574\"),\n    @(\"42\u00d723=966\", \"63\u00d729=1827\"),\n    @(\"65\u00d766=4290\", \"41\u00d791=3731\"),\n    @(\"89\u00d786=7654\", \"81\u00d752=4212\"),\n    @(\"75\u00d793=6975\", \"28\u00d775=2100\"),\n    @(\"18\u00d758=1044\", \"97\u00d772=6984\"),\n    @(\"54\u00d728=1512\", \"51\u00d781=4131\"),\n    @(\"95\u00d758=5510\", \"44\u00d793=4092\"),\n    @(\"15\u00d711=165\", \"92\u00d732=2944\"),\n    @(\"40\u00d712=480\", \"81\u00d751=4131\"),\n    @(\"41\u00d764=2624\", \"71\u00d783=5893\"),\n    @(\"57\u00d722=1254\", \"46\u00d790=4140\"),\n    @(\"81\u00d794=7614\", \"82\u00d771=5822\"),\n    @(\"47\u00d766=3102\", \"70\u00d768=4760\"),\n    @(\"89\u00d749=4361\", \"28\u00d725=700\"),\n    @(\"80\u00d738=3040\", \"68\u00d721=1428\"),\n    @(\"68\u00d760=4080\", \"39\u00d716=624\"),\n    @(\"94\u00d738=3572\", \"55\u00d793=5115\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $found = $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n\n    if (-not $found) {\n        Write-Output \"WARNING: text not found for replacement: $oldText\"\n    }\n}\n\nWrite-Output \"Done\"\n"}
